$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Triangle test-case data rows (row 2 .. row 20), columns:
#   A = case no. (sequential), B/C/D = side lengths, E = expected, F = actual, G = result
$rows = @(
    @(0,   100, 100, "长度超出范围", "长度超出范围", "pass"),
    @(1,   100, 100, "等腰三角形",   "等腰三角形",   "pass"),
    @(10,  100, 95,  "普通三角形",   "普通三角形",   "pass"),
    @(100, 100, 100, "等边三角形",   "等边三角形",   "pass"),
    @(190, 100, 90,  "不构成三角形", "不构成三角形", "pass"),
    @(200, 100, 110, "普通三角形",   "普通三角形",   "pass"),
    @(201, 100, 100, "长度超出范围", "长度超出范围", "pass"),
    @(100, 0,   100, "长度超出范围", "长度超出范围", "pass"),
    @(100, 1,   90,  "不构成三角形", "不构成三角形", "pass"),
    @(100, 10,  105, "普通三角形",   "普通三角形",   "pass"),
    @(100, 190, 110, "普通三角形",   "普通三角形",   "pass"),
    @(100, 200, 100, "不构成三角形", "不构成三角形", "pass"),
    @(100, 201, 100, "长度超出范围", "长度超出范围", "pass"),
    @(100, 100, 0,   "长度超出范围", "长度超出范围", "pass"),
    @(100, 110, 1,   "不构成三角形", "不构成三角形", "pass"),
    @(100, 100, 10,  "等腰三角形",   "等腰三角形",   "pass"),
    @(100, 95,  190, "普通三角形",   "普通三角形",   "pass"),
    @(100, 90,  200, "不构成三角形", "不构成三角形", "pass"),
    @(100, 110, 201, "长度超出范围", "长度超出范围", "pass")
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $r - 1
    $ws.Cells.Item($r, 2).Value = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $row[3]
    $ws.Cells.Item($r, 6).Value = $row[4]
    $ws.Cells.Item($r, 7).Value = $row[5]
    $r = $r + 1
}

$ws.Range("E22").Select()
